$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Change the date heading "06.03.2023" -> "27.03.2023"
#    (this is the first / only occurrence at this point in the document)
# ---------------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute("06.03.2023", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "27.03.2023", 2)

# ---------------------------------------------------------------------------
# 2) Insert three new paragraphs right before the bullet
#    "Fehlerkorrektur: Einstellungen der Anlagen 11, 17 und 18 ..."
#    which is paragraph #4 after the rename above.
#
#    New content (in order):
#      a) bullet: "Fehlerkorrektur: In Anlage 5b wurden die Vornamen nicht
#                  gemäß den Einstellungen in der INI-Datei übernommen."
#      b) bullet: "Anlagen 10, 11, 17 und 18 verwenden nun ebenfalls die
#                  Vornamenseinstellungen aus der INI-Datei statt die der
#                  Serienbriefvorlagen (Ausnahme Adressfeld)."
#      c) bold heading (no list): "06.03.2023"
# ---------------------------------------------------------------------------
$anchorPara = $d.Paragraphs.Item(4)
$anchorRange = $anchorPara.Range
$anchorRange.InsertParagraphBefore()
$anchorRange.InsertParagraphBefore()
$anchorRange.InsertParagraphBefore()

# a) first new bullet - inherits the "Listenabsatz" / numId 5 list style already
$bullet1 = $d.Paragraphs.Item(4)
$bullet1.Range.Text = "Fehlerkorrektur: In Anlage 5b wurden die Vornamen nicht gemäß den Einstellungen in der INI-Datei übernommen."

# b) second new bullet - also inherits the correct list style
$bullet2 = $d.Paragraphs.Item(5)
$bullet2.Range.Text = "Anlagen 10, 11, 17 und 18 verwenden nun ebenfalls die Vornamenseinstellungen aus der INI-Datei statt die der Serienbriefvorlagen (Ausnahme Adressfeld)."

# c) new bold date heading - remove the inherited list formatting/style first
$dateHeading = $d.Paragraphs.Item(6)
$dateHeading.Range.ListFormat.RemoveNumbers()
$dateHeading.Style = "Standard"
$dateHeading.Range.Text = "06.03.2023"
$dateHeading.Range.Font.Bold = 1
$dateHeading.Range.Font.BoldBi = 1
